$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header (row 15): merged banner, same look as the row 7 banner ---
$ws.Range("A7:E7").Copy()
$ws.Range("B15:F15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B15").Value = "try to sort 2 or 3 columns not entire column"
$ws.Range("B15:F15").Merge()

# --- New table header (row 16): same look as the row 1 header ---
$ws.Range("A1:C1").Copy()
$ws.Range("A16:C16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(16).RowHeight = 23
$ws.Range("A16").Value = "Department"
$ws.Range("B16").Value = "Salary"

# --- Sort the Department/Salary columns (B2:C5) by Salary ascending, spilling into A17:B20 ---
$ws.Range("A17:B20").FormulaArray = "=SORT(B2:C5,2,1)"

# --- Match the author's final selection state ---
$ws.Range("C17").Select()

$wb.Save()
